$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C48").Value = "[name=`"Talulah`"]  My 'teacher' wasn’t anything like you, for sure.`n"
$ws.Range("C70").Value = "[name=`"Talulah`"]  You already know all the rest. After FrostNova and I had our 'friendly discussion,' our group successfully joined up with the guerrillas.`n"
$ws.Range("C83").Value = "[name=`"Talulah`"]  A friend of mine in a certain city calls that kind of behavior 'reunion.'`n"
$ws.Range("C84").Value = "[name=`"Talulah`"]  He calls on the Infected to come together, naming it the 'Reunion Movement,' revolving around one identical belief, protesting against Ursus’s cruel rule over the Infected.`n"
$ws.Range("C87").Value = "[name=`"FrostNova`"]  If you want to go south to 'protest,' we could just save ourselves the trouble and find a division to kill us.`n"
$ws.Range("C90").Value = "[name=`"Talulah`"]  'You are not alone.' `n"
$ws.Range("C123").Value = "[name=`"FrostNova`"]  What I’m thinking is written on my face. I don’t know anything you’d call a 'future.'`n"
$ws.Range("D41").Value = "[name=`"탈룰라`"]  나도 프로스트노바가 어렸을 때 어떻게 자랐는지는 잘 몰라. 교양이 있는 걸로 봐선, 그래도 나름 잘 살았던 것 같던데.`n"
$ws.Range("D42").Value = "[name=`"탈룰라`"]  다른 아이들이 내가 겪었던 일을 겪게 놔두고 싶진 않아.`n"
$ws.Range("D45").Value = "[name=`"알리나`"]  ……탈룰라, 이 땅이 아이들에게 어떤 경험을 선사할지, 우리보다 더 터무니없는 경험을 하게 할지는, 아무도 알 수 없는 거야.`n"
